$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M5").Value = 1.07
$ws.Range("N5").Value = 9
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 3
$ws.Range("Q5").Value = 2.15
$ws.Range("R5").Value = 1.67
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.08
$ws.Range("I8").Value = 8
$ws.Range("K8").Value = 2.4
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("W8").Value = 7
$ws.Range("X8").Value = 6.5
$ws.Range("Z8").Value = 9
$ws.Range("AA8").Value = 12
$ws.Range("AB8").Value = 29
$ws.Range("AC8").Value = 11
$ws.Range("AD8").Value = 9
$ws.Range("AJ8").Value = 81
$ws.Range("AN8").Value = 3.25
$ws.Range("AO8").Value = 6.5
$ws.Range("AP8").Value = 19
$ws.Range("AQ8").Value = 19
$ws.Range("AU8").Value = 9.5
$ws.Range("AZ8").Value = 151
$ws.Range("G9").Value = 3.5
$ws.Range("I9").Value = 2.25
$ws.Range("L9").Value = 3
$ws.Range("O9").Value = 1.36
$ws.Range("P9").Value = 3.2
$ws.Range("Q9").Value = 2.15
$ws.Range("R9").Value = 1.67
$ws.Range("U9").Value = 1.91
$ws.Range("V9").Value = 1.91
$ws.Range("AM9").Value = 301
$ws.Range("AO9").Value = 19
$ws.Range("AV9").Value = 51
$ws.Range("M10").Value = 1.1
$ws.Range("N10").Value = 7
$ws.Range("I14").Value = 1.55
$ws.Range("R14").Value = 1.8
$ws.Range("G18").Value = 7.7
$ws.Range("H18").Value = 4.9
$ws.Range("I18").Value = 1.34
$ws.Range("K18").Value = 2.75
$ws.Range("L18").Value = 1.72
$ws.Range("P18").Value = 6.1
$ws.Range("Q18").Value = 1.31
$ws.Range("R18").Value = 3.15
$ws.Range("T18").Value = 4.15
$ws.Range("Y18").Value = 24
$ws.Range("AC18").Value = 10.5
$ws.Range("AD18").Value = 11.5
$ws.Range("AE18").Value = 14.5
$ws.Range("AG18").Value = 12.5
$ws.Range("AI18").Value = 8.75
$ws.Range("AN18").Value = 9.75
$ws.Range("AT18").Value = 4.15
$ws.Range("AU18").Value = 6.8
$ws.Range("AX18").Value = 6.1
$ws.Range("AZ18").Value = 14
$ws.Range("BA18").Value = 26
